$d = $word.ActiveDocument
$r = $d.Content
$r.Collapse(0)

$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:spacing w:line="240" w:lineRule="auto"/>
    <w:rPr>
      <w:lang w:val="en-US"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:t>Eventually getting</w:t>
  </w:r>
</w:p>
<w:tbl>
  <w:tblPr>
    <w:tblStyle w:val="TableGrid"/>
    <w:tblW w:w="0" w:type="auto"/>
    <w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/>
  </w:tblPr>
  <w:tblGrid>
    <w:gridCol w:w="4814"/>
    <w:gridCol w:w="4814"/>
  </w:tblGrid>
  <w:tr>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="4814" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:rPr>
            <w:lang w:val="en-US"/>
          </w:rPr>
        </w:pPr>
        <w:r>
          <w:rPr>
            <w:lang w:val="en-US"/>
          </w:rPr>
          <w:t xml:space="preserve">FO :  ****** Current epoch: 30 ****** </w:t>
        </w:r>
      </w:p>
      <w:p>
        <w:pPr>
          <w:rPr>
            <w:lang w:val="en-US"/>
          </w:rPr>
        </w:pPr>
        <w:r>
          <w:rPr>
            <w:lang w:val="en-US"/>
          </w:rPr>
          <w:t xml:space="preserve">INFO : Training sample: </w:t>
        </w:r>
        <w:r>
          <w:rPr>
            <w:lang w:val="en-US"/>
          </w:rPr>
          <w:tab/>
          <w:t xml:space="preserve"> 60/ 300 ...</w:t>
        </w:r>
      </w:p>
      <w:p>
        <w:pPr>
          <w:rPr>
            <w:lang w:val="en-US"/>
          </w:rPr>
        </w:pPr>
        <w:r>
          <w:rPr>
            <w:lang w:val="en-US"/>
          </w:rPr>
          <w:t xml:space="preserve">INFO : Training sample: </w:t>
        </w:r>
        <w:r>
          <w:rPr>
            <w:lang w:val="en-US"/>
          </w:rPr>
          <w:tab/>
          <w:t xml:space="preserve"> 120/ 300 ...</w:t>
        </w:r>
      </w:p>
      <w:p>
        <w:pPr>
          <w:rPr>
            <w:lang w:val="en-US"/>
          </w:rPr>
        </w:pPr>
        <w:r>
          <w:rPr>
            <w:lang w:val="en-US"/>
          </w:rPr>
          <w:t xml:space="preserve">INFO : Training sample: </w:t>
        </w:r>
        <w:r>
          <w:rPr>
            <w:lang w:val="en-US"/>
          </w:rPr>
          <w:tab/>
          <w:t xml:space="preserve"> 180/ 300 ...</w:t>
        </w:r>
      </w:p>
      <w:p>
        <w:pPr>
          <w:rPr>
            <w:lang w:val="en-US"/>
          </w:rPr>
        </w:pPr>
        <w:r>
          <w:rPr>
            <w:lang w:val="en-US"/>
          </w:rPr>
          <w:t xml:space="preserve">INFO : Training sample: </w:t>
        </w:r>
        <w:r>
          <w:rPr>
            <w:lang w:val="en-US"/>
          </w:rPr>
          <w:tab/>
          <w:t xml:space="preserve"> 240/ 300 ...</w:t>
        </w:r>
      </w:p>
      <w:p>
        <w:pPr>
          <w:rPr>
            <w:lang w:val="en-US"/>
          </w:rPr>
        </w:pPr>
        <w:r>
          <w:rPr>
            <w:lang w:val="en-US"/>
          </w:rPr>
          <w:t xml:space="preserve">INFO : Training sample: </w:t>
        </w:r>
        <w:r>
          <w:rPr>
            <w:lang w:val="en-US"/>
          </w:rPr>
          <w:tab/>
          <w:t xml:space="preserve"> 300/ 300 ...</w:t>
        </w:r>
      </w:p>
      <w:p>
        <w:pPr>
          <w:rPr>
            <w:lang w:val="en-US"/>
          </w:rPr>
        </w:pPr>
        <w:r>
          <w:rPr>
            <w:lang w:val="en-US"/>
          </w:rPr>
          <w:t>INFO : loss=3.99 ; accuracy=0.973</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="4814" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:rPr>
            <w:lang w:val="en-US"/>
          </w:rPr>
        </w:pPr>
        <w:r>
          <w:rPr>
            <w:lang w:val="en-US"/>
          </w:rPr>
          <w:t>INFO : Evaluation</w:t>
        </w:r>
      </w:p>
      <w:p>
        <w:pPr>
          <w:rPr>
            <w:lang w:val="en-US"/>
          </w:rPr>
        </w:pPr>
        <w:r>
          <w:rPr>
            <w:lang w:val="en-US"/>
          </w:rPr>
          <w:t xml:space="preserve">INFO : Sample: </w:t>
        </w:r>
        <w:r>
          <w:rPr>
            <w:lang w:val="en-US"/>
          </w:rPr>
          <w:tab/>
          <w:t xml:space="preserve"> 30/ 150 ...</w:t>
        </w:r>
      </w:p>
      <w:p>
        <w:pPr>
          <w:rPr>
            <w:lang w:val="en-US"/>
          </w:rPr>
        </w:pPr>
        <w:r>
          <w:rPr>
            <w:lang w:val="en-US"/>
          </w:rPr>
          <w:t xml:space="preserve">INFO : Sample: </w:t>
        </w:r>
        <w:r>
          <w:rPr>
            <w:lang w:val="en-US"/>
          </w:rPr>
          <w:tab/>
          <w:t xml:space="preserve"> 60/ 150 ...</w:t>
        </w:r>
      </w:p>
      <w:p>
        <w:pPr>
          <w:rPr>
            <w:lang w:val="en-US"/>
          </w:rPr>
        </w:pPr>
        <w:r>
          <w:rPr>
            <w:lang w:val="en-US"/>
          </w:rPr>
          <w:t xml:space="preserve">INFO : Sample: </w:t>
        </w:r>
        <w:r>
          <w:rPr>
            <w:lang w:val="en-US"/>
          </w:rPr>
          <w:tab/>
          <w:t xml:space="preserve"> 90/ 150 ...</w:t>
        </w:r>
      </w:p>
      <w:p>
        <w:pPr>
          <w:rPr>
            <w:lang w:val="en-US"/>
          </w:rPr>
        </w:pPr>
        <w:r>
          <w:rPr>
            <w:lang w:val="en-US"/>
          </w:rPr>
          <w:t xml:space="preserve">INFO : Sample: </w:t>
        </w:r>
        <w:r>
          <w:rPr>
            <w:lang w:val="en-US"/>
          </w:rPr>
          <w:tab/>
          <w:t xml:space="preserve"> 120/ 150 ...</w:t>
        </w:r>
      </w:p>
      <w:p>
        <w:pPr>
          <w:rPr>
            <w:lang w:val="en-US"/>
          </w:rPr>
        </w:pPr>
        <w:r>
          <w:rPr>
            <w:lang w:val="en-US"/>
          </w:rPr>
          <w:t xml:space="preserve">INFO : Sample: </w:t>
        </w:r>
        <w:r>
          <w:rPr>
            <w:lang w:val="en-US"/>
          </w:rPr>
          <w:tab/>
          <w:t xml:space="preserve"> 150/ 150 ...</w:t>
        </w:r>
      </w:p>
      <w:p>
        <w:pPr>
          <w:rPr>
            <w:lang w:val="en-US"/>
          </w:rPr>
        </w:pPr>
        <w:r>
          <w:rPr>
            <w:lang w:val="en-US"/>
          </w:rPr>
          <w:t>INFO : loss=4.26 ; accuracy=0.76</w:t>
        </w:r>
      </w:p>
      <w:p>
        <w:pPr>
          <w:rPr>
            <w:lang w:val="en-US"/>
          </w:rPr>
        </w:pPr>
        <w:r>
          <w:rPr>
            <w:lang w:val="en-US"/>
          </w:rPr>
          <w:t xml:space="preserve">INFO : </w:t>
        </w:r>
        <w:proofErr w:type="spellStart"/>
        <w:r>
          <w:rPr>
            <w:lang w:val="en-US"/>
          </w:rPr>
          <w:t>validation_loss</w:t>
        </w:r>
        <w:proofErr w:type="spellEnd"/>
        <w:r>
          <w:rPr>
            <w:lang w:val="en-US"/>
          </w:rPr>
          <w:t xml:space="preserve">=4.259 ; </w:t>
        </w:r>
        <w:proofErr w:type="spellStart"/>
        <w:r>
          <w:rPr>
            <w:lang w:val="en-US"/>
          </w:rPr>
          <w:t>best_validation_loss</w:t>
        </w:r>
        <w:proofErr w:type="spellEnd"/>
        <w:r>
          <w:rPr>
            <w:lang w:val="en-US"/>
          </w:rPr>
          <w:t>=4.302</w:t>
        </w:r>
      </w:p>
    </w:tc>
  </w:tr>
</w:tbl>
<w:p>
  <w:pPr>
    <w:spacing w:line="240" w:lineRule="auto"/>
    <w:rPr>
      <w:lang w:val="en-US"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p>
  <w:pPr>
    <w:spacing w:line="240" w:lineRule="auto"/>
    <w:rPr>
      <w:lang w:val="en-US"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:t>Fourth step: train the model on the entirety of the training set (15k instances), checking the entirety of the validation set (3k instances)</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:t>. Remember to leave the saved model where it is.</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:br/>
    <w:t>(and to write the code to load a pre-defined model)</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:spacing w:line="240" w:lineRule="auto"/>
    <w:rPr>
      <w:lang w:val="en-US"/>
    </w:rPr>
  </w:pPr>
</w:p>

'@

$r.InsertXML($xml)
